$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to text format first so numeric-looking values
# (e.g. "102.43") are stored as text, matching the original inlineStr cells,
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "51.083.90"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").Value = "2.959.95"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").Value = "380.12"
$ws.Range("E5").Value = "  +1.08%  "

# Row 6
$ws.Range("D6").Value = "102.43"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value = "  +1.76%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "0.588"
$ws.Range("E9").Value = "  +0.66%  "

# Row 10
$ws.Range("D10").Value = "36.53"
$ws.Range("E10").Value = "  -0.40%  "

# Row 11
$ws.Range("E11").Value = "  -0.46%  "

# Row 12
$ws.Range("D12").Value = "0.0852"
$ws.Range("E12").Value = "  +2.04%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.425.00"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
$ws.Range("D14").Value = "18.41"
$ws.Range("E14").Value = "  +2.55%  "

# Row 15
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").Value = "12.39"
$ws.Range("E15").Value = "  +73.87%  "

# Row 16
$ws.Range("D16").Value = "7.75"
$ws.Range("E16").Value = "  +5.44%  "

# Row 17
$ws.Range("D17").Value = "2.960.44"
$ws.Range("E17").Value = "  +0.80%  "

# Row 18
$ws.Range("D18").Value = "1.01"
$ws.Range("E18").Value = "  +3.83%  "

# Row 19
$ws.Range("D19").Value = "51.136.85"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("E20").Value = "  -2.88%  "

# Row 21
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -0.86%  "

# Row 22
$ws.Range("E22").Value = "  +1.00%  "

# Row 23
$ws.Range("E23").Value = "  +18.02%  "

# Row 24
$ws.Range("D24").Value = "268.98"
$ws.Range("E24").Value = "  +2.25%  "

# Row 25
$ws.Range("D25").Value = "69.79"

# Row 26
$ws.Range("D26").Value = "8.01"
$ws.Range("E26").Value = "  -2.08%  "

# Row 27
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("D28").Value = "0.166"
$ws.Range("E28").Value = "  -0.38%  "

# Row 29
$ws.Range("D29").Value = "25.91"
$ws.Range("E29").Value = "  +0.96%  "

# Row 30
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -9.51%  "

# Row 31
$ws.Range("E31").Value = "  -3.75%  "

# Row 32
$ws.Range("D32").Value = "10.54"
$ws.Range("E32").Value = "  +7.12%  "

# Row 33
$ws.Range("D33").Value = "51.14"
$ws.Range("E33").Value = "  +0.95%  "

# Row 34
$ws.Range("D34").Value = "34.15"
$ws.Range("E34").Value = "  +0.20%  "

# Row 35
$ws.Range("E35").Value = "  +2.10%  "

# Row 36
$ws.Range("D36").Value = "0.0436"
$ws.Range("E36").Value = "  -4.11%  "

# Row 37
$ws.Range("E37").Value = "  +0.15%  "

# Row 38
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  +9.36%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.117"
$ws.Range("E39").Value = "  +2.09%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "16.72"
$ws.Range("E40").Value = "  +1.54%  "

# Row 41
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  +2.76%  "

# Row 42
$ws.Range("D42").Value = "2.49"
$ws.Range("E42").Value = "  -3.69%  "

# Row 43
$ws.Range("D43").Value = "123.90"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "21.78"
$ws.Range("E44").Value = "  +3.07%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.58"
$ws.Range("E45").Value = "  +10.86%  "

# Row 46
$ws.Range("D46").Value = "2.088.74"
$ws.Range("E46").Value = "  +4.07%  "

# Row 47
$ws.Range("E47").Value = "  -1.01%  "

# Row 48
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -0.63%  "

# Row 50
$ws.Range("D50").Value = "0.0323"
$ws.Range("E50").Value = "  -7.11%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "5.35"
$ws.Range("E51").Value = "  +6.71%  "

# Restore the default cell style for the Price column so no extra
# number-format styling is introduced (keeps cells visually "General").
$ws.Range("D2:D51").Style = "Normal"